$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a serial date value of 45179 (2023-09-10)
# for every data row (rows 2 through 420). Update it to 45180 (2023-09-11).
$ws.Range("C2:C420").Value = 45180
